# Update the "Sprint 4" burndown-chart source data.
# Root-cause edits (everything else -- the B4:B17 "ideal" curve, the
# C4:C10 "actual" curve and the chart caches -- recalculates from these):
#   - B3 (and C3, which mirrors it) drops from 101 to 94
#   - F10 (effort logged against the 7th task) goes from 0 to 7, which
#     makes C10 pick up a value for the first time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")

$ws.Range("B3").Value = 94
$ws.Range("C3").Value = 94
$ws.Range("F10").Value = 7

# Re-enter the "Actual" formula across C4:C10 as a single block fill so it
# extends onto the newly-logged C10 (Excel turns this into one shared
# formula group spanning C4:C10, matching a fill-down/fill-handle edit).
$ws.Range("C4:C10").Formula = "=C3-F4"

$excel.Calculate()

$ws.Range("C15").Select()
